# Bulk changes and improvements for loading footings file and asserting equal.
# Update the "mapping" column (C) on the hidden "__footings__" sheet:
# replace bracket-style mapping references like "[parameter.a]" with
# slash-style references like "/parameter.a/".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__footings__")

$rows = @(39, 40, 41, 42, 43, 44, 45, 46, 60, 61, 75, 76, 90, 91, 93, 94)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    if ($current -ne $null -and $current.StartsWith("[") -and $current.EndsWith("]")) {
        $inner = $current.Substring(1, $current.Length - 2)
        $cell.Value = "/" + $inner + "/"
    }
}
